$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.20950566666667
$ws.Range("H2").Value = 87.628517
$ws.Range("I2").Value = 0.01829497698069002
$ws.Range("J2").Value = 0.01840828041918582
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.06706333333333332
$ws.Range("N2").Value = 0.20119
$ws.Range("O2").Value = 0.01564803973866319
$ws.Range("P2").Value = 0.01878890697104963
$ws.Range("Q2").Value = 1.958886815025555
$ws.Range("R2").Value = 17.62998133523
$ws.Range("S2").Value = 0.0002862805268117657
$ws.Range("T2").Value = 0.0003458714682930769
$ws.Range("G3").Value = 29.20950566666667
$ws.Range("H3").Value = 87.628517
$ws.Range("I3").Value = 0.01829497698069002
$ws.Range("J3").Value = 0.01840828041918582
$ws.Range("O3").Value = 0.2726197454399388
$ws.Range("P3").Value = 0.3273398534952746
$ws.Range("Q3").Value = 34.12767565629544
$ws.Range("R3").Value = 307.149080906659
$ws.Range("S3").Value = 0.004987571967305254
$ws.Range("T3").Value = 0.006025763815516219
$ws.Range("G4").Value = 29.20950566666667
$ws.Range("H4").Value = 87.628517
$ws.Range("I4").Value = 0.01829497698069002
$ws.Range("J4").Value = 0.01840828041918582
$ws.Range("M4").Value = 0.72155
$ws.Range("N4").Value = 2.16465
$ws.Range("O4").Value = 0.1683608987539007
$ws.Range("P4").Value = 0.2021542197668005
$ws.Range("Q4").Value = 21.07611881378334
$ws.Range("R4").Value = 189.68506932405
$ws.Range("S4").Value = 0.003080158767150896
$ws.Range("T4").Value = 0.003721311565388981
$ws.Range("G5").Value = 29.20950566666667
$ws.Range("H5").Value = 87.628517
$ws.Range("I5").Value = 0.01829497698069002
$ws.Range("J5").Value = 0.01840828041918582
$ws.Range("M5").Value = 2.1492875
$ws.Range("N5").Value = 4.298575
$ws.Range("O5").Value = 0.501498129277977
$ws.Range("P5").Value = 0.4014390664699025
$ws.Range("Q5").Value = 62.77962541054583
$ws.Range("R5").Value = 376.677752463275
$ws.Range("S5").Value = 0.009174896730999697
$ws.Range("T5").Value = 0.007389802906794141
$ws.Range("G6").Value = 29.20950566666667
$ws.Range("H6").Value = 87.628517
$ws.Range("I6").Value = 0.01829497698069002
$ws.Range("J6").Value = 0.01840828041918582
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1794573333333333
$ws.Range("N6").Value = 0.538372
$ws.Range("O6").Value = 0.04187318678952025
$ws.Range("P6").Value = 0.05027795329697268
$ws.Range("Q6").Value = 5.241859994924889
$ws.Range("R6").Value = 47.176739954324
$ws.Range("S6").Value = 0.0007660689884224065
$ws.Range("T6").Value = 0.0009255306631934014
$ws.Range("I7").Value = 0.913374480506715
$ws.Range("J7").Value = 0.9190311407684336
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.06706333333333332
$ws.Range("N7").Value = 0.20119
$ws.Range("O7").Value = 0.01564803973866319
$ws.Range("P7").Value = 0.01878890697104963
$ws.Range("Q7").Value = 97.79718383542553
$ws.Range("R7").Value = 880.1746545188299
$ws.Range("S7").Value = 0.01429252016724992
$ws.Range("T7").Value = 0.01726759060739572
$ws.Range("I8").Value = 0.913374480506715
$ws.Range("J8").Value = 0.9190311407684336
$ws.Range("O8").Value = 0.2726197454399388
$ws.Range("P8").Value = 0.3273398534952746
$ws.Range("S8").Value = 0.249003918367077
$ws.Range("T8").Value = 0.3008355189767342
$ws.Range("I9").Value = 0.913374480506715
$ws.Range("J9").Value = 0.9190311407684336
$ws.Range("M9").Value = 0.72155
$ws.Range("N9").Value = 2.16465
$ws.Range("O9").Value = 0.1683608987539007
$ws.Range("P9").Value = 0.2021542197668005
$ws.Range("Q9").Value = 1052.222645207783
$ws.Range("R9").Value = 9470.003806870049
$ws.Range("S9").Value = 0.1537765484369877
$ws.Range("T9").Value = 0.1857860232034353
$ws.Range("I10").Value = 0.913374480506715
$ws.Range("J10").Value = 0.9190311407684336
$ws.Range("M10").Value = 2.1492875
$ws.Range("N10").Value = 4.298575
$ws.Range("O10").Value = 0.501498129277977
$ws.Range("P10").Value = 0.4014390664699025
$ws.Range("Q10").Value = 3134.265093981045
$ws.Range("R10").Value = 18805.59056388627
$ws.Range("S10").Value = 0.4580555933043616
$ws.Range("T10").Value = 0.3689350032068495
$ws.Range("I11").Value = 0.913374480506715
$ws.Range("J11").Value = 0.9190311407684336
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1794573333333333
$ws.Range("N11").Value = 0.538372
$ws.Range("O11").Value = 0.04187318678952025
$ws.Range("P11").Value = 0.05027795329697268
$ws.Range("Q11").Value = 261.6992169384449
$ws.Range("R11").Value = 2355.292952446004
$ws.Range("S11").Value = 0.0382459002310387
$ws.Range("T11").Value = 0.04620700477401883
$ws.Range("G12").Value = 57.98602933333333
$ws.Range("H12").Value = 173.958088
$ws.Range("I12").Value = 0.03631876156896331
$ws.Range("J12").Value = 0.03654368891224535
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.06706333333333332
$ws.Range("N12").Value = 0.20119
$ws.Range("O12").Value = 0.01564803973866319
$ws.Range("P12").Value = 0.01878890697104963
$ws.Range("Q12").Value = 3.888736413857777
$ws.Range("R12").Value = 34.99862772472
$ws.Range("S12").Value = 0.0005683174242901713
$ws.Range("T12").Value = 0.0006866159713511558
$ws.Range("G13").Value = 57.98602933333333
$ws.Range("H13").Value = 173.958088
$ws.Range("I13").Value = 0.03631876156896331
$ws.Range("J13").Value = 0.03654368891224535
$ws.Range("O13").Value = 0.2726197454399388
$ws.Range("P13").Value = 0.3273398534952746
$ws.Range("Q13").Value = 67.74946567968621
$ws.Range("R13").Value = 609.7451911171759
$ws.Range("S13").Value = 0.009901211533624612
$ws.Range("T13").Value = 0.01196220577471128
$ws.Range("G14").Value = 57.98602933333333
$ws.Range("H14").Value = 173.958088
$ws.Range("I14").Value = 0.03631876156896331
$ws.Range("J14").Value = 0.03654368891224535
$ws.Range("M14").Value = 0.72155
$ws.Range("N14").Value = 2.16465
$ws.Range("O14").Value = 0.1683608987539007
$ws.Range("P14").Value = 0.2021542197668005
$ws.Range("Q14").Value = 41.83981946546666
$ws.Range("R14").Value = 376.5583751891999
$ws.Range("S14").Value = 0.006114659339379291
$ws.Range("T14").Value = 0.007387460919455637
$ws.Range("G15").Value = 57.98602933333333
$ws.Range("H15").Value = 173.958088
$ws.Range("I15").Value = 0.03631876156896331
$ws.Range("J15").Value = 0.03654368891224535
$ws.Range("M15").Value = 2.1492875
$ws.Range("N15").Value = 4.298575
$ws.Range("O15").Value = 0.501498129277977
$ws.Range("P15").Value = 0.4014390664699025
$ws.Range("Q15").Value = 124.6286480207666
$ws.Range("R15").Value = 747.7718881245999
$ws.Range("S15").Value = 0.01821379098452798
$ws.Range("T15").Value = 0.0146700643622983
$ws.Range("G16").Value = 57.98602933333333
$ws.Range("H16").Value = 173.958088
$ws.Range("I16").Value = 0.03631876156896331
$ws.Range("J16").Value = 0.03654368891224535
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1794573333333333
$ws.Range("N16").Value = 0.538372
$ws.Range("O16").Value = 0.04187318678952025
$ws.Range("P16").Value = 0.05027795329697268
$ws.Range("Q16").Value = 10.40601819474844
$ws.Range("R16").Value = 93.65416375273598
$ws.Range("S16").Value = 0.00152078228714125
$ws.Range("T16").Value = 0.00183734188442897
$ws.Range("G17").Value = 29.481085
$ws.Range("H17").Value = 58.96217
$ws.Range("I17").Value = 0.01846507700595112
$ws.Range("J17").Value = 0.01238628926567028
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.06706333333333332
$ws.Range("N17").Value = 0.20119
$ws.Range("O17").Value = 0.01564803973866319
$ws.Range("P17").Value = 0.01878890697104963
$ws.Range("Q17").Value = 1.977099830383333
$ws.Range("R17").Value = 11.8625989823
$ws.Range("S17").Value = 0.0002889422587665991
$ws.Range("T17").Value = 0.0002327248367291895
$ws.Range("G18").Value = 29.481085
$ws.Range("H18").Value = 58.96217
$ws.Range("I18").Value = 0.01846507700595112
$ws.Range("J18").Value = 0.01238628926567028
$ws.Range("O18").Value = 0.2726197454399388
$ws.Range("P18").Value = 0.3273398534952746
$ws.Range("Q18").Value = 34.44498234093167
$ws.Range("R18").Value = 206.66989404559
$ws.Range("S18").Value = 0.005033944592891263
$ws.Range("T18").Value = 0.004054526113574601
$ws.Range("G19").Value = 29.481085
$ws.Range("H19").Value = 58.96217
$ws.Range("I19").Value = 0.01846507700595112
$ws.Range("J19").Value = 0.01238628926567028
$ws.Range("M19").Value = 0.72155
$ws.Range("N19").Value = 2.16465
$ws.Range("O19").Value = 0.1683608987539007
$ws.Range("P19").Value = 0.2021542197668005
$ws.Range("Q19").Value = 21.27207688175
$ws.Range("R19").Value = 127.6324612905
$ws.Range("S19").Value = 0.003108796960281917
$ws.Range("T19").Value = 0.002503940642307471
$ws.Range("G20").Value = 29.481085
$ws.Range("H20").Value = 58.96217
$ws.Range("I20").Value = 0.01846507700595112
$ws.Range("J20").Value = 0.01238628926567028
$ws.Range("M20").Value = 2.1492875
$ws.Range("N20").Value = 4.298575
$ws.Range("O20").Value = 0.501498129277977
$ws.Range("P20").Value = 0.4014390664699025
$ws.Range("Q20").Value = 63.36332747693749
$ws.Range("R20").Value = 253.45330990775
$ws.Range("S20").Value = 0.009260201575458277
$ws.Range("T20").Value = 0.00497234039983685
$ws.Range("G21").Value = 29.481085
$ws.Range("H21").Value = 58.96217
$ws.Range("I21").Value = 0.01846507700595112
$ws.Range("J21").Value = 0.01238628926567028
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.1794573333333333
$ws.Range("N21").Value = 0.538372
$ws.Range("O21").Value = 0.04187318678952025
$ws.Range("P21").Value = 0.05027795329697268
$ws.Range("Q21").Value = 5.290596897873333
$ws.Range("R21").Value = 31.74358138724
$ws.Range("S21").Value = 0.0007731916185530668
$ws.Range("T21").Value = 0.0006227572732221642
$ws.Range("G22").Value = 21.628479
$ws.Range("H22").Value = 64.885437
$ws.Range("I22").Value = 0.01354670393768061
$ws.Range("J22").Value = 0.01363060063446486
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.06706333333333332
$ws.Range("N22").Value = 0.20119
$ws.Range("O22").Value = 0.01564803973866319
$ws.Range("P22").Value = 0.01878890697104963
$ws.Range("Q22").Value = 1.45047789667
$ws.Range("R22").Value = 13.05430107003
$ws.Range("S22").Value = 0.0002119793615447313
$ws.Range("T22").Value = 0.0002561040872804904
$ws.Range("G23").Value = 21.628479
$ws.Range("H23").Value = 64.885437
$ws.Range("I23").Value = 0.01354670393768061
$ws.Range("J23").Value = 0.01363060063446486
$ws.Range("O23").Value = 0.2726197454399388
$ws.Range("P23").Value = 0.3273398534952746
$ws.Range("Q23").Value = 25.270188570611
$ws.Range("R23").Value = 227.431697135499
$ws.Range("S23").Value = 0.003693098979040705
$ws.Range("T23").Value = 0.004461838814738326
$ws.Range("G24").Value = 21.628479
$ws.Range("H24").Value = 64.885437
$ws.Range("I24").Value = 0.01354670393768061
$ws.Range("J24").Value = 0.01363060063446486
$ws.Range("M24").Value = 0.72155
$ws.Range("N24").Value = 2.16465
$ws.Range("O24").Value = 0.1683608987539007
$ws.Range("P24").Value = 0.2021542197668005
$ws.Range("Q24").Value = 15.60602902245
$ws.Range("R24").Value = 140.45426120205
$ws.Range("S24").Value = 0.002280735250100913
$ws.Range("T24").Value = 0.0027554834362131
$ws.Range("G25").Value = 21.628479
$ws.Range("H25").Value = 64.885437
$ws.Range("I25").Value = 0.01354670393768061
$ws.Range("J25").Value = 0.01363060063446486
$ws.Range("M25").Value = 2.1492875
$ws.Range("N25").Value = 4.298575
$ws.Range("O25").Value = 0.501498129277977
$ws.Range("P25").Value = 0.4014390664699025
$ws.Range("Q25").Value = 46.48581955871249
$ws.Range("R25").Value = 278.914917352275
$ws.Range("S25").Value = 0.006793646682629431
$ws.Range("T25").Value = 0.005471855594123634
$ws.Range("G26").Value = 21.628479
$ws.Range("H26").Value = 64.885437
$ws.Range("I26").Value = 0.01354670393768061
$ws.Range("J26").Value = 0.01363060063446486
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.1794573333333333
$ws.Range("N26").Value = 0.538372
$ws.Range("O26").Value = 0.04187318678952025
$ws.Range("P26").Value = 0.05027795329697268
$ws.Range("Q26").Value = 3.881389165396
$ws.Range("R26").Value = 34.932502488564
$ws.Range("S26").Value = 0.0005672436643648298
$ws.Range("T26").Value = 0.0006853187021093106
